$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: delete B2, D2, E2; update C2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 24.82869911336519
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: update B3:E3
$ws.Range("B3").Value = 22.112663767150138
$ws.Range("C3").Value = 31.332085099570122
$ws.Range("D3").Value = 33.035550857034877
$ws.Range("E3").Value = 16.342677613145952

# Update selection to B1:E3
$ws.Range("B1:E3").Select()
